# ModelComponentClassDiagram.pptx — "Updated images in dev guide"
#
# The AddressBook/Person model-component class diagram on slide 1 was
# reworked into a TaskManager/Task diagram:
#   - AddressBook         -> TaskManager
#   - UniquePersonList    -> UniqueTaskList
#   - Person              -> Task
#   - ReadOnlyPerson      -> ReadOnlyTask            (2nd line of the
#                                                      <<interface>> box)
#   - ReadOnlyAddressBook -> ReadOnlyTaskManager      (2nd line of the
#                                                      <<interface>> box)
#   - the Person attribute box "Name"  -> "Description" (widened)
#   - the Person attribute box "Phone" -> "DateTime"     (widened)
#   - the Person attribute box "Email" -> "TaskStatus"   (widened)
#   - the Person attribute box "Address" (and its connector) was removed
#     entirely, since Task only keeps 3 attributes

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

# --- simple single-run text renames ---------------------------------
(Get-ShapeById $s 46).TextFrame.TextRange.Text  = "TaskManager"       # AddressBook
(Get-ShapeById $s 49).TextFrame.TextRange.Text  = "UniqueTaskList"    # UniquePersonList
(Get-ShapeById $s 62).TextFrame.TextRange.Text  = "Task"              # Person

# --- two-line "<<interface>>\n<Name>" boxes: only touch the 2nd line -
$sh72 = Get-ShapeById $s 72
$tr72 = $sh72.TextFrame.TextRange
$len72 = $tr72.Text.Length
$tr72.Characters(15, $len72 - 14).Text = "ReadOnlyTask"        # ReadOnlyPerson

$sh100 = Get-ShapeById $s 100
$tr100 = $sh100.TextFrame.TextRange
$len100 = $tr100.Text.Length
$tr100.Characters(15, $len100 - 14).Text = "ReadOnlyTaskManager"  # ReadOnlyAddressBook

# --- Person attribute boxes: rename + widen --------------------------
$sh76 = Get-ShapeById $s 76
$sh76.TextFrame.TextRange.Text = "Description"                 # Name
$sh76.Width = 822003.5 / 12700.0
$sh76.Left  = $sh76.Left

$sh80 = Get-ShapeById $s 80
$sh80.TextFrame.TextRange.Text = "DateTime"                    # Phone
$sh80.Width = 822002.5 / 12700.0

$sh83 = Get-ShapeById $s 83
$sh83.TextFrame.TextRange.Text = "TaskStatus"                  # Email
$sh83.Width = 822002.5 / 12700.0

# --- Task no longer has an "Address" attribute: drop the box plus its
#     connector to the decision diamond (delete highest id first so
#     indices of shapes we still need stay stable) -------------------
(Get-ShapeById $s 86).Delete()   # Elbow Connector 85
(Get-ShapeById $s 85).Delete()   # Rectangle "Address"
